$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990207675/1749682509529/zqlwK29jaz-camera-video-segment-driver-1749682514529.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOBNLZN6S%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T190024Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJGMEQCIEKdza4yTvRWKD6rkdBl5Py4hhXkoxWP%2FW52iirnYIsTAiBYcnUN6z6ufs52j%2BVvfGe0RJ5D%2Fu4tYBSJ9caDNsRYySrmAwjy%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMtf4svLDQjFDbCixbKroD3aqjYKQFyOKc06nIazkhqdX8VZCPiNsYBoSo6gRawdgP4x8pj3yfbexcrm%2FJRFKpmKqgvdSga%2FnHLimKxza%2FpfqdBfzVBdOB2azjZoQpRHJ8dr188Yb96MR%2BZxQyw1lW%2BgffUYo9HmM9zIIaBfpYZm0jrZNgzZJkj1%2BBjLWTy5ami%2F%2BAvyoJ5f5uPIJAQWxea0mwGuYyy1OZD21BYo7HaSFWo4nT0XSQjeZcOxQEmldfFOaMDi%2B7LYwFTl0a%2FvghlnBSD%2BVZ4B0k9C1UST6hj%2Bd7oRc4xCKFa37GDVt2Z3X9tQ%2BF1Ph7XPpV5gJb0Suj7EAawPv42fPIfX%2F%2BYXoHV2peQqdRuIY26wmEMHrh3E1Pvb72JLHrmYVWvDjPeKiPki4irC0VyNvc8BElgLfyPn9vwaASpYh4iYd9UUBFVxOrGe3ISKiZ7uKpGP3L8ElVSrFXQ6l7kyEzBA4Ov3C385MMzP3jQSjLIzlqc%2BqmPk%2FkIGtmd7H4Nbk%2FCeXR0bBkFnWKd9lBIc8Hlpw7czzFNDZfelvRkWboT5W8g8XPXZ1EMKLFNBJnNsrHcH9rV1hfBNp37e59dtIDpTC3kazCBjqmAWUroY9XQqZcC4oHq0c9D26n%2BWrTej4Dxzl6osqtphbzqIpm6%2BCGMd%2F4FT64VJVa502fv1pcIG7uBGug7HgnOchFrS4ovEfZRl9VQzhaj0P38MbwvthFD9PLUdw9ew%2FcascHf1bMR1rTqzNMdAIyM8mXycEaCohP2fM%2FJmx0iN2kX5XK3aGtMeb7RX6aJBMXlkYtjlquT8um1kNHCAnkyQbbr0jQhEI%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2003%3A00%3A24%20GMT&X-Amz-Signature=9859d1a7208929048cad01b2bbb2b8ec295ed5391a5612fd52d72e83ee7e23fd'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990207675/1749682509529/JfEAwSUnRy-camera-video-segment-1749682514529.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOBNLZN6S%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T190024Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJGMEQCIEKdza4yTvRWKD6rkdBl5Py4hhXkoxWP%2FW52iirnYIsTAiBYcnUN6z6ufs52j%2BVvfGe0RJ5D%2Fu4tYBSJ9caDNsRYySrmAwjy%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMtf4svLDQjFDbCixbKroD3aqjYKQFyOKc06nIazkhqdX8VZCPiNsYBoSo6gRawdgP4x8pj3yfbexcrm%2FJRFKpmKqgvdSga%2FnHLimKxza%2FpfqdBfzVBdOB2azjZoQpRHJ8dr188Yb96MR%2BZxQyw1lW%2BgffUYo9HmM9zIIaBfpYZm0jrZNgzZJkj1%2BBjLWTy5ami%2F%2BAvyoJ5f5uPIJAQWxea0mwGuYyy1OZD21BYo7HaSFWo4nT0XSQjeZcOxQEmldfFOaMDi%2B7LYwFTl0a%2FvghlnBSD%2BVZ4B0k9C1UST6hj%2Bd7oRc4xCKFa37GDVt2Z3X9tQ%2BF1Ph7XPpV5gJb0Suj7EAawPv42fPIfX%2F%2BYXoHV2peQqdRuIY26wmEMHrh3E1Pvb72JLHrmYVWvDjPeKiPki4irC0VyNvc8BElgLfyPn9vwaASpYh4iYd9UUBFVxOrGe3ISKiZ7uKpGP3L8ElVSrFXQ6l7kyEzBA4Ov3C385MMzP3jQSjLIzlqc%2BqmPk%2FkIGtmd7H4Nbk%2FCeXR0bBkFnWKd9lBIc8Hlpw7czzFNDZfelvRkWboT5W8g8XPXZ1EMKLFNBJnNsrHcH9rV1hfBNp37e59dtIDpTC3kazCBjqmAWUroY9XQqZcC4oHq0c9D26n%2BWrTej4Dxzl6osqtphbzqIpm6%2BCGMd%2F4FT64VJVa502fv1pcIG7uBGug7HgnOchFrS4ovEfZRl9VQzhaj0P38MbwvthFD9PLUdw9ew%2FcascHf1bMR1rTqzNMdAIyM8mXycEaCohP2fM%2FJmx0iN2kX5XK3aGtMeb7RX6aJBMXlkYtjlquT8um1kNHCAnkyQbbr0jQhEI%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2003%3A00%3A24%20GMT&X-Amz-Signature=817db515998ce2970d90e9fee039c25be72acd2ab21eb77c79cfe7223ac09d8a'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991154589/1749681707753/UtIzP4FPML-camera-video-segment-driver-1749681712753.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOBNLZN6S%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T190024Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJGMEQCIEKdza4yTvRWKD6rkdBl5Py4hhXkoxWP%2FW52iirnYIsTAiBYcnUN6z6ufs52j%2BVvfGe0RJ5D%2Fu4tYBSJ9caDNsRYySrmAwjy%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMtf4svLDQjFDbCixbKroD3aqjYKQFyOKc06nIazkhqdX8VZCPiNsYBoSo6gRawdgP4x8pj3yfbexcrm%2FJRFKpmKqgvdSga%2FnHLimKxza%2FpfqdBfzVBdOB2azjZoQpRHJ8dr188Yb96MR%2BZxQyw1lW%2BgffUYo9HmM9zIIaBfpYZm0jrZNgzZJkj1%2BBjLWTy5ami%2F%2BAvyoJ5f5uPIJAQWxea0mwGuYyy1OZD21BYo7HaSFWo4nT0XSQjeZcOxQEmldfFOaMDi%2B7LYwFTl0a%2FvghlnBSD%2BVZ4B0k9C1UST6hj%2Bd7oRc4xCKFa37GDVt2Z3X9tQ%2BF1Ph7XPpV5gJb0Suj7EAawPv42fPIfX%2F%2BYXoHV2peQqdRuIY26wmEMHrh3E1Pvb72JLHrmYVWvDjPeKiPki4irC0VyNvc8BElgLfyPn9vwaASpYh4iYd9UUBFVxOrGe3ISKiZ7uKpGP3L8ElVSrFXQ6l7kyEzBA4Ov3C385MMzP3jQSjLIzlqc%2BqmPk%2FkIGtmd7H4Nbk%2FCeXR0bBkFnWKd9lBIc8Hlpw7czzFNDZfelvRkWboT5W8g8XPXZ1EMKLFNBJnNsrHcH9rV1hfBNp37e59dtIDpTC3kazCBjqmAWUroY9XQqZcC4oHq0c9D26n%2BWrTej4Dxzl6osqtphbzqIpm6%2BCGMd%2F4FT64VJVa502fv1pcIG7uBGug7HgnOchFrS4ovEfZRl9VQzhaj0P38MbwvthFD9PLUdw9ew%2FcascHf1bMR1rTqzNMdAIyM8mXycEaCohP2fM%2FJmx0iN2kX5XK3aGtMeb7RX6aJBMXlkYtjlquT8um1kNHCAnkyQbbr0jQhEI%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2003%3A00%3A24%20GMT&X-Amz-Signature=2dbb175659cf1c17ff3738dd4b6e011b6c37408c53499f157b9e49c60f2da840'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991154589/1749681707753/418ikFdSWa-camera-video-segment-1749681712753.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSOBNLZN6S%2F20250612%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250612T190024Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEBkaCXVzLXdlc3QtMiJGMEQCIEKdza4yTvRWKD6rkdBl5Py4hhXkoxWP%2FW52iirnYIsTAiBYcnUN6z6ufs52j%2BVvfGe0RJ5D%2Fu4tYBSJ9caDNsRYySrmAwjy%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F8BEAQaDDc4MTIwNDk0MjI0NCIMtf4svLDQjFDbCixbKroD3aqjYKQFyOKc06nIazkhqdX8VZCPiNsYBoSo6gRawdgP4x8pj3yfbexcrm%2FJRFKpmKqgvdSga%2FnHLimKxza%2FpfqdBfzVBdOB2azjZoQpRHJ8dr188Yb96MR%2BZxQyw1lW%2BgffUYo9HmM9zIIaBfpYZm0jrZNgzZJkj1%2BBjLWTy5ami%2F%2BAvyoJ5f5uPIJAQWxea0mwGuYyy1OZD21BYo7HaSFWo4nT0XSQjeZcOxQEmldfFOaMDi%2B7LYwFTl0a%2FvghlnBSD%2BVZ4B0k9C1UST6hj%2Bd7oRc4xCKFa37GDVt2Z3X9tQ%2BF1Ph7XPpV5gJb0Suj7EAawPv42fPIfX%2F%2BYXoHV2peQqdRuIY26wmEMHrh3E1Pvb72JLHrmYVWvDjPeKiPki4irC0VyNvc8BElgLfyPn9vwaASpYh4iYd9UUBFVxOrGe3ISKiZ7uKpGP3L8ElVSrFXQ6l7kyEzBA4Ov3C385MMzP3jQSjLIzlqc%2BqmPk%2FkIGtmd7H4Nbk%2FCeXR0bBkFnWKd9lBIc8Hlpw7czzFNDZfelvRkWboT5W8g8XPXZ1EMKLFNBJnNsrHcH9rV1hfBNp37e59dtIDpTC3kazCBjqmAWUroY9XQqZcC4oHq0c9D26n%2BWrTej4Dxzl6osqtphbzqIpm6%2BCGMd%2F4FT64VJVa502fv1pcIG7uBGug7HgnOchFrS4ovEfZRl9VQzhaj0P38MbwvthFD9PLUdw9ew%2FcascHf1bMR1rTqzNMdAIyM8mXycEaCohP2fM%2FJmx0iN2kX5XK3aGtMeb7RX6aJBMXlkYtjlquT8um1kNHCAnkyQbbr0jQhEI%3D&X-Amz-SignedHeaders=host&response-expires=Fri%2C%2013%20Jun%202025%2003%3A00%3A24%20GMT&X-Amz-Signature=686dfb73c12554bf1cdd73b9bd2886c4078ff298b88d83f714269ab773d119d2'
